{"js": "// The edit described by the diff boils down to a single visible-text\n// change: a new sentence is appended to the very end of the document's\n// last paragraph (\"...the score of the game. \" -> \"...the score of the\n// game. The scoring function is only used by the AI. \"). Every other\n// hunk in the diff only merges/re-splits <w:r> runs and drops/repositions\n// spell/grammar-check <w:proofErr/> markers around unchanged text, which\n// is not visible content and is not something the Word/Office.js object\n// model exposes a way to control directly (Word itself regenerates those\n// markers from its proofing engine) - so there is nothing else to apply.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"The scoring function is only used by the AI. \", \"End\");\nawait context.sync();\n", "ps1": "# The edit described by the diff boils down to a single visible-text\n# change: a new sentence is appended to the very end of the document's\n# last paragraph (\"...the score of the game. \" -> \"...the score of the\n# game. The scoring function is only used by the AI. \"). Every other\n# hunk in the diff only merges/re-splits runs and drops/repositions\n# spell/grammar-check proofErr markers around otherwise-unchanged text,\n# which is not visible content and isn't something the Word object model\n# exposes a way to control directly (Word regenerates those markers from\n# its own proofing engine) - so there is nothing else to apply.\n\n$d = $word.ActiveDocument\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertAfter(\"The scoring function is only used by the AI. \")\n"}
